# Added a code snippet to convert R file to Rmd

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 5 data ------------------------------------------------------
# Shared-strings must be created in this order (title, description, tags,
# url, src) to match the target sharedStrings.xml ordering.
$ws.Range("A5").Value = 'Convert an `R` file into an Rmd'
$ws.Range("D5").Value = 'Use knitr::spin() to convert R file into Rmd'
$ws.Range("E5").Value = 'R; Convert to Rmd'
$ws.Range("C5").Value = 'https://github.com/sciencificity/convert-r-to-rmd'
$ws.Range("B5").Value = 'images/markus-spiske-hGb5WqRrWIg-unsplash.jpg'
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1

# --- Hyperlink on C5 -------------------------------------------------------
$h = $ws.Hyperlinks.Add($ws.Range("C5"), 'https://github.com/sciencificity/convert-r-to-rmd')
$ws.Range("C5").Style = "Hyperlink"

# --- Column width adjustments (best-effort match to target widths) --------
$ws.Columns.Item(1).ColumnWidth = 27.666666666666668
$ws.Columns.Item(2).ColumnWidth = 33.666666666666664
$ws.Columns.Item(3).ColumnWidth = 57.666666666666664
$ws.Columns.Item(4).ColumnWidth = 22.333333333333332
$ws.Columns.Item(5).ColumnWidth = 28.666666666666668

# --- Selection update -------------------------------------------------------
$ws.Range("B5").Select() | Out-Null

# --- Enable iterative calculation (workbook setting) ------------------------
$excel.Iteration = $true
$excel.MaxIterations = 100
$excel.MaxChange = 0.001
